# Added ACH Payment Entry Bootstrap Page Reusable Library
# Extend the ACHData sheet with a new "TypeOfAccount" column (Saving/Checking)
# and a third sample data row.

$wb  = $excel.ActiveWorkbook
$ach = $wb.Worksheets.Item("ACHData")
$udf = $wb.Worksheets.Item("UDFData")

# New header column E: TypeOfAccount
$ach.Range("E1").Value = "TypeOfAccount"
$ach.Range("E2").Value = "Saving"
$ach.Range("E3").Value = "Checking"

# New sample row (id 3) mirroring the row-2 account/routing values, account type Checking
$ach.Range("A4").Value = 3
$ach.Range("B4").Value = 95125489
$ach.Range("C4").Value = 95125489
$ach.Range("D4").Value = 256072691
$ach.Range("E4").Value = "Checking"

# Move the active selection: UDFData loses focus (now D7), ACHData becomes the
# active/selected sheet with C9 selected.
$udf.Range("D7").Select()

$ach.Activate()
$ach.Range("C9").Select()
